$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "SamplesTab" row (row 3) query had a typo (`sample_id_id`) in the
# ORDER BY clause plus the LIMIT clause jammed on the same line. Fix it to
# match the corrected query used elsewhere (and reflow LIMIT onto its own
# line), exactly as committed upstream.
$fixedQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
  WHERE diag.disease_term IN ['Glioma']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
Order By samp.sample_id ASC 
LIMIT 100
'@

$ws.Range("B3").Value = $fixedQuery

# Update the active selection left after editing (Excel drops the prior
# multi-cell / scrolled selection once the edit is committed).
$ws.Range("B4").Select()
